$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the data range keeps text formatting so numeric-looking values
# (e.g. prices, ranks) are stored as text, matching the source data.
$ws.Range("B2:G51").NumberFormat = "@"

$ws.Range('D2').Value = '245.76'
$ws.Range('G2').Value = '10'
$ws.Range('D3').Value = '26.35'
$ws.Range('G3').Value = '10'
$ws.Range('D4').Value = '5.132'
$ws.Range('G4').Value = '10'
$ws.Range('G5').Value = '10'
$ws.Range('G6').Value = '10'
$ws.Range('D7').Value = '3.022'
$ws.Range('G7').Value = '10'
$ws.Range('D8').Value = '0.8147'
$ws.Range('G8').Value = '10'
$ws.Range('D9').Value = '0.8482'
$ws.Range('G9').Value = '10'
$ws.Range('D10').Value = '0.1338'
$ws.Range('G10').Value = '10'
$ws.Range('B11').Value = 'BitrueCoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D11').Value = '0.02851'
$ws.Range('E11').Value = '10BitrueCoinBTR'
$ws.Range('G11').Value = '10'
$ws.Range('B12').Value = 'BitMartToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D12').Value = '0.09374'
$ws.Range('E12').Value = '11BitMartTokenBMX'
$ws.Range('G12').Value = '10'
$ws.Range('B13').Value = 'BitForexToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D13').Value = '0.001516'
$ws.Range('E13').Value = '12BitForexTokenBF'
$ws.Range('G13').Value = '10'
$ws.Range('B14').Value = 'TigerCash'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D14').Value = '0.006239'
$ws.Range('E14').Value = '13TigerCashTCH'
$ws.Range('G14').Value = '10'
$ws.Range('B15').Value = 'LEO'
$ws.Range('C15').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D15').Value = '3.552'
$ws.Range('E15').Value = '14LEOLEO'
$ws.Range('G15').Value = '10'
$ws.Range('B16').Value = 'BTSEToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D16').Value = '2.118'
$ws.Range('E16').Value = '15BTSETokenBTSE'
$ws.Range('G16').Value = '10'
$ws.Range('B17').Value = 'One'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D17').Value = '0.0006019'
$ws.Range('E17').Value = '16OneONE'
$ws.Range('G17').Value = '10'
$ws.Range('B18').Value = 'BitpandaEcosystemToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D18').Value = '0.3134'
$ws.Range('E18').Value = '17BitpandaEcosystemTokenBEST'
$ws.Range('G18').Value = '10'
$ws.Range('B19').Value = 'MandalaExchangeToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D19').Value = '0.06949'
$ws.Range('E19').Value = '18MandalaExchangeTokenMDX'
$ws.Range('G19').Value = '10'
$ws.Range('B20').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C20').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D20').Value = '0.03221'
$ws.Range('E20').Value = '19LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('G20').Value = '10'
$ws.Range('G21').Value = '10'
$ws.Range('D22').Value = '3.748'
$ws.Range('G22').Value = '10'
$ws.Range('D23').Value = '0.04687'
$ws.Range('G23').Value = '10'
$ws.Range('G24').Value = '10'
$ws.Range('D25').Value = '0.001247'
$ws.Range('G25').Value = '10'
$ws.Range('G26').Value = '10'
$ws.Range('D27').Value = '0.00009599'
$ws.Range('G27').Value = '10'
$ws.Range('E28').Value = '27UpBotsUBXTWorstin24h'
$ws.Range('G28').Value = '10'
$ws.Range('G29').Value = '10'
$ws.Range('G30').Value = '10'
$ws.Range('G31').Value = '10'
$ws.Range('G32').Value = '10'
$ws.Range('G33').Value = '10'
$ws.Range('G34').Value = '10'
$ws.Range('G35').Value = '10'
$ws.Range('G36').Value = '10'
$ws.Range('G37').Value = '10'
$ws.Range('G38').Value = '10'
$ws.Range('G39').Value = '10'
$ws.Range('D40').Value = '0.03650'
$ws.Range('G40').Value = '10'
$ws.Range('B41').Value = 'KickToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D41').Value = '0.006110'
$ws.Range('E41').Value = '40KickTokenKICKBestin24h'
$ws.Range('G41').Value = '10'
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D42').Value = '0.1054'
$ws.Range('E42').Value = '41BKEXTokenBKK'
$ws.Range('G42').Value = '10'
$ws.Range('D43').Value = '0.002564'
$ws.Range('G43').Value = '10'
$ws.Range('D44').Value = '0.008110'
$ws.Range('G44').Value = '10'
$ws.Range('D45').Value = '0.00005306'
$ws.Range('G45').Value = '10'
$ws.Range('G46').Value = '10'
$ws.Range('D47').Value = '0.1723'
$ws.Range('E47').Value = '46CoinbaseStockTokenCOIN'
$ws.Range('G47').Value = '10'
$ws.Range('D48').Value = '0.002048'
$ws.Range('G48').Value = '10'
$ws.Range('G49').Value = '10'
$ws.Range('G50').Value = '10'
$ws.Range('G51').Value = '10'
